$d = $word.ActiveDocument

# --- Change 1: add "Prof. Dr. " before "Heng Ji" in the Strong-styled
#     speaker-name paragraph (NOT the Heading3 "Heng Ji" higher up). ---
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    $txt = $para.Range.Text.Trim()
    if ($txt -eq "Heng Ji" -and $para.Style.NameLocal -eq "Normal (Web)") {
        $target = $para
    }
}

if ($target -ne $null) {
    $r = $d.Range($target.Range.Start, $target.Range.End)
    $r.Find.Execute("Heng Ji", $false, $false, $false, $false, $false, $true,
                     0, $false, "Prof. Dr. Heng Ji", 1)
}

# --- Change 2: rewrite the bio paragraph so the misspell-checked
#     "Woese" run merges back with its neighbours (drops the
#     w:proofErr spell-check wrapper runs introduced by authoring). ---
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -match "^is a Tenured Full Professor") {
        $bioRange = $d.Range($para.Range.Start, $para.Range.End)
        $needle = "Carl R. Woese Institute"
        $bioRange.Find.Execute($needle, $false, $false, $false, $false, $false,
                                $true, 0, $false, $needle, 2)
    }
}
